$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 985
$ws.Range("I16").Value = 985
$ws.Range("K16").Value = 985
$ws.Range("M16").Value = -755
$ws.Range("H19").Value = 557.6429000000001
$ws.Range("I19").Value = 465.6875
$ws.Range("J19").Value = 680.25
$ws.Range("K19").Value = 465.6875
$ws.Range("L19").Value = 680.25
$ws.Range("M19").Value = -290.6875
$ws.Range("N19").Value = -1030.25
$ws.Range("H38").Value = 470.9091
$ws.Range("I38").Value = 168
$ws.Range("J38").Value = 3500
$ws.Range("K38").Value = 504
$ws.Range("L38").Value = 10500
$ws.Range("M38").Value = -132
$ws.Range("N38").Value = -11244
$ws.Range("H40").Value = 4103.75
$ws.Range("J40").Value = 4477.0835
$ws.Range("L40").Value = 4477.0835
$ws.Range("N40").Value = -4827.0835
$ws.Range("H45").Value = 2061.2856
$ws.Range("I45").Value = 789
$ws.Range("J45").Value = 2159.1538
$ws.Range("K45").Value = 2367
$ws.Range("L45").Value = 6477.4614
$ws.Range("M45").Value = -2175
$ws.Range("N45").Value = -6861.4614
$ws.Range("H98").Value = 189208
$ws.Range("I98").Value = 1810.52
$ws.Range("J98").Value = 858484.7
$ws.Range("K98").Value = 1810.52
$ws.Range("L98").Value = 858484.7
$ws.Range("M98").Value = -312.52
$ws.Range("N98").Value = -861480.7
$ws.Range("H116").Value = 16100
$ws.Range("I116").Value = 14082.1
$ws.Range("K116").Value = 14082.1
$ws.Range("M116").Value = -10640.1
$ws.Range("H122").Value = 189208
$ws.Range("I122").Value = 1810.52
$ws.Range("J122").Value = 858484.7
$ws.Range("K122").Value = 5431.559999999999
$ws.Range("L122").Value = 2575454.1
$ws.Range("M122").Value = -2981.559999999999
$ws.Range("N122").Value = -2580354.1
$ws.Range("H127").Value = 1145.5
$ws.Range("H132").Value = 1817.65
$ws.Range("I132").Value = 1820.7297
$ws.Range("K132").Value = 5462.189100000001
$ws.Range("M132").Value = -2932.189100000001
$ws.Range("H134").Value = 71849.336
$ws.Range("J134").Value = 71849.336
$ws.Range("L134").Value = 71849.336
$ws.Range("N134").Value = -81989.336
$ws.Range("H138").Value = 2401.9114
$ws.Range("I138").Value = 1183.3125
$ws.Range("K138").Value = 3549.9375
$ws.Range("M138").Value = 1590.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25645774
$ws.Range("I32").Value = 25645774
$ws.Range("K32").Value = 25645774
$ws.Range("M32").Value = -25645487
$ws.Range("H45").Value = 4214.15
$ws.Range("I45").Value = 1915.75
$ws.Range("J45").Value = 7661.75
$ws.Range("K45").Value = 1915.75
$ws.Range("L45").Value = 7661.75
$ws.Range("M45").Value = -1538.75
$ws.Range("N45").Value = -8415.75
$ws.Range("H74").Value = 9526807
$ws.Range("I74").Value = 10419046
$ws.Range("K74").Value = 10419046
$ws.Range("M74").Value = -10418172
$ws.Range("H77").Value = 9526807
$ws.Range("I77").Value = 10419046
$ws.Range("K77").Value = 52095230
$ws.Range("M77").Value = -52090862
$ws.Range("H110").Value = 2327.7856
$ws.Range("I110").Value = 1060.6364
$ws.Range("J110").Value = 6974
$ws.Range("K110").Value = 1060.6364
$ws.Range("L110").Value = 6974
$ws.Range("M110").Value = 984.3635999999999
$ws.Range("N110").Value = -11064
$ws.Range("H114").Value = 32000
$ws.Range("J114").Value = 32000
$ws.Range("L114").Value = 32000
$ws.Range("N114").Value = -40678
$ws.Range("H132").Value = 4992.3057
$ws.Range("I132").Value = 3993.375
$ws.Range("K132").Value = 11980.125
$ws.Range("M132").Value = -9450.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3365.077
$ws.Range("I94").Value = 1704.1818
$ws.Range("K94").Value = 1704.1818
$ws.Range("M94").Value = -1253.1818
$ws.Range("H134").Value = 1522.1562
$ws.Range("I134").Value = 1003.8
$ws.Range("J134").Value = 3373.4285
$ws.Range("K134").Value = 3011.4
$ws.Range("L134").Value = 10120.2855
$ws.Range("M134").Value = -476.3999999999996
$ws.Range("N134").Value = -15190.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23703.963
$ws.Range("I31").Value = 2513.7273
$ws.Range("K31").Value = 2513.7273
$ws.Range("M31").Value = -2218.7273
$ws.Range("H34").Value = 23703.963
$ws.Range("I34").Value = 2513.7273
$ws.Range("K34").Value = 2513.7273
$ws.Range("M34").Value = -2311.7273
$ws.Range("H99").Value = 2985.75
$ws.Range("J99").Value = 3188.4
$ws.Range("L99").Value = 3188.4
$ws.Range("N99").Value = -6184.4
$ws.Range("H122").Value = 7737.8335
$ws.Range("I122").Value = 2748.3333
$ws.Range("J122").Value = 12727.333
$ws.Range("K122").Value = 8244.999899999999
$ws.Range("L122").Value = 38181.999
$ws.Range("M122").Value = -5794.999899999999
$ws.Range("N122").Value = -43081.999
$ws.Range("H126").Value = 2985.75
$ws.Range("J126").Value = 3188.4
$ws.Range("L126").Value = 9565.200000000001
$ws.Range("N126").Value = -14505.2
$ws.Range("H132").Value = 2219
$ws.Range("I132").Value = 1928.2046
$ws.Range("K132").Value = 5784.6138
$ws.Range("M132").Value = -3254.6138
$ws.Range("H141").Value = 233217.5
$ws.Range("J141").Value = 233217.5
$ws.Range("L141").Value = 233217.5
$ws.Range("N141").Value = -243577.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1760
$ws.Range("I11").Value = 2942.5
$ws.Range("J11").Value = 183.33333
$ws.Range("K11").Value = 8827.5
$ws.Range("L11").Value = 549.99999
$ws.Range("M11").Value = -8687.5
$ws.Range("N11").Value = -829.99999
$ws.Range("H26").Value = 929.3
$ws.Range("I26").Value = 927.4167
$ws.Range("J26").Value = 932.125
$ws.Range("K26").Value = 2782.2501
$ws.Range("L26").Value = 2796.375
$ws.Range("M26").Value = -2494.2501
$ws.Range("N26").Value = -3372.375
$ws.Range("H80").Value = 10996.167
$ws.Range("I80").Value = 9997.5
$ws.Range("K80").Value = 29992.5
$ws.Range("M80").Value = -29056.5
$ws.Range("H83").Value = 10996.167
$ws.Range("I83").Value = 9997.5
$ws.Range("K83").Value = 89977.5
$ws.Range("M83").Value = -85297.5
$ws.Range("H107").Value = 90911400
$ws.Range("I107").Value = 2376.1667
$ws.Range("K107").Value = 7128.500100000001
$ws.Range("M107").Value = -5208.500100000001
$ws.Range("H120").Value = 13454.6
$ws.Range("I120").Value = 4909.2
$ws.Range("K120").Value = 14727.6
$ws.Range("M120").Value = -9889.599999999999
$ws.Range("H132").Value = 4522.276
$ws.Range("I132").Value = 4164.4287
$ws.Range("J132").Value = 4856.2666
$ws.Range("K132").Value = 37479.85830000001
$ws.Range("L132").Value = 43706.3994
$ws.Range("M132").Value = -34949.85830000001
$ws.Range("N132").Value = -48766.3994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3327.7273
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3560.5
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 3560.5
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -5556.5
$ws.Range("H83").Value = 3327.7273
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3560.5
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 17802.5
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -27786.5
$ws.Range("H102").Value = 3510723.5
$ws.Range("I102").Value = 4763657
$ws.Range("J102").Value = 2509.6
$ws.Range("K102").Value = 4763657
$ws.Range("L102").Value = 2509.6
$ws.Range("M102").Value = -4762035
$ws.Range("N102").Value = -5753.6
$ws.Range("H128").Value = 48110.25
$ws.Range("J128").Value = 48110.25
$ws.Range("L128").Value = 48110.25
$ws.Range("N128").Value = -58070.25
$ws.Range("H139").Value = 74897
$ws.Range("J139").Value = 74897
$ws.Range("L139").Value = 74897
$ws.Range("N139").Value = -85177

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9020.105
$ws.Range("I7").Value = 5697.273
$ws.Range("J7").Value = 13589
$ws.Range("K7").Value = 5697.273
$ws.Range("L7").Value = 13589
$ws.Range("M7").Value = -5585.273
$ws.Range("N7").Value = -13813
$ws.Range("H16").Value = 432.13043
$ws.Range("I16").Value = 432.13043
$ws.Range("K16").Value = 432.13043
$ws.Range("M16").Value = -262.13043
$ws.Range("H22").Value = 1518.25
$ws.Range("I22").Value = 813.05554
$ws.Range("J22").Value = 2787.6
$ws.Range("K22").Value = 813.05554
$ws.Range("L22").Value = 2787.6
$ws.Range("M22").Value = -518.05554
$ws.Range("N22").Value = -3377.6
$ws.Range("H27").Value = 1518.25
$ws.Range("I27").Value = 813.05554
$ws.Range("J27").Value = 2787.6
$ws.Range("K27").Value = 813.05554
$ws.Range("L27").Value = 2787.6
$ws.Range("M27").Value = -706.05554
$ws.Range("N27").Value = -3001.6
$ws.Range("H93").Value = 1362.5555
$ws.Range("I93").Value = 1262.1666
$ws.Range("K93").Value = 1262.1666
$ws.Range("M93").Value = -14.16660000000002
$ws.Range("H126").Value = 9020.105
$ws.Range("I126").Value = 5697.273
$ws.Range("J126").Value = 13589
$ws.Range("K126").Value = 17091.819
$ws.Range("L126").Value = 40767
$ws.Range("M126").Value = -14621.819
$ws.Range("N126").Value = -45707
$ws.Range("H132").Value = 2861.2952
$ws.Range("I132").Value = 2428.5386
$ws.Range("J132").Value = 5361.6665
$ws.Range("K132").Value = 7285.6158
$ws.Range("L132").Value = 16084.9995
$ws.Range("M132").Value = -4755.6158
$ws.Range("N132").Value = -21144.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10150
$ws.Range("I122").Value = 11798.333
$ws.Range("K122").Value = 35394.999
$ws.Range("M122").Value = -32944.999
$ws.Range("H137").Value = 65879
$ws.Range("J137").Value = 65879
$ws.Range("L137").Value = 65879
$ws.Range("N137").Value = -76079
